# Apply data corrections to the "Inscricoes" sheet (Table1)
# This mirrors an upstream data refresh: several numeric counts
# (Inscritos / Pagos / Inscrições homologadas columns) were updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$updates = @{
    "E2" = 77
    "F2" = 53
    "H2" = 65
    "E3" = 31
    "F3" = 19
    "H3" = 23
    "E4" = 31
    "E5" = 8
    "F5" = 4
    "H5" = 4
    "E6" = 12
    "F6" = 5
    "H6" = 6
    "E9" = 36
    "F9" = 19
    "H9" = 28
    "E15" = 197
    "F15" = 117
    "H15" = 158
    "E16" = 14
    "E17" = 152
    "F17" = 81
    "H17" = 113
    "E18" = 150
    "E19" = 86
    "F19" = 50
    "H19" = 63
    "E20" = 9
    "E21" = 4
    "E22" = 10
    "F22" = 7
    "H22" = 7
    "E23" = 10
    "E24" = 31
    "F24" = 20
    "H24" = 24
    "E26" = 43
    "F26" = 25
    "H26" = 35
    "E27" = 23
    "E29" = 21
    "E31" = 3
    "E33" = 52
    "E36" = 138
    "F36" = 69
    "H36" = 101
    "E37" = 73
    "F37" = 44
    "H37" = 56
    "E38" = 101
    "F38" = 28
    "H38" = 48
    "E39" = 32
    "E40" = 31
    "F40" = 22
    "H40" = 24
    "E41" = 52
    "F41" = 28
    "H41" = 39
    "E42" = 48
    "F42" = 33
    "H42" = 42
    "E43" = 39
    "F43" = 23
    "H43" = 26
    "E44" = 37
    "F44" = 23
    "H44" = 33
    "F45" = 15
    "H45" = 22
    "E46" = 35
    "F46" = 14
    "H46" = 23
    "E47" = 77
    "F47" = 51
    "H47" = 61
    "E48" = 53
    "F48" = 35
    "H48" = 41
    "F49" = 49
    "H49" = 66
    "E50" = 34
    "E57" = 22
    "E60" = 25
    "E61" = 40
    "E62" = 64
    "E63" = 54
    "E64" = 45
    "F64" = 23
    "H64" = 28
    "F65" = 17
    "H65" = 30
    "E66" = 43
    "E67" = 49
    "F67" = 28
    "H67" = 36
    "E69" = 20
    "E70" = 57
    "F70" = 31
    "H70" = 44
    "E71" = 54
    "E72" = 57
    "E73" = 39
    "E74" = 23
    "E76" = 63
    "E77" = 69
    "F77" = 29
    "H77" = 46
    "E78" = 53
    "F78" = 24
    "H78" = 45
    "E79" = 52
    "F79" = 27
    "H79" = 38
    "E81" = 24
    "F81" = 16
    "H81" = 21
    "E86" = 6
    "E87" = 23
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
